$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new "Senior Researcher, Paul C. Lauterbur Lab at SIAT"
# entry (with its bullet sub-points and two trailing blank paragraphs) at the
# very top of the document body, before the existing first paragraph.
# ---------------------------------------------------------------------------
$newContentXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="859"/>
      <w:tab w:val="left" w:pos="860"/>
      <w:tab w:val="right" w:pos="11080"/>
    </w:tabs>
    <w:ind w:left="144"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="4472C4" w:themeColor="accent1"/>
    </w:rPr>
    <w:t>Senior Researcher, Paul C. Lauterbur Lab at SIAT</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Shenzhen, CN </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">| </w:t>
  </w:r>
  <w:r>
    <w:t>Nov 2016 - Jan 2017</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="859"/>
      <w:tab w:val="left" w:pos="860"/>
      <w:tab w:val="right" w:pos="11080"/>
    </w:tabs>
    <w:ind w:left="144"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:i/>
    </w:rPr>
    <w:t>EMG signal pattern recognition for hand gestures using spectral analysis</w:t>
  </w:r>
  <w:r>
    <w:tab/>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="859"/>
      <w:tab w:val="left" w:pos="860"/>
    </w:tabs>
    <w:spacing w:line="244" w:lineRule="exact"/>
  </w:pPr>
  <w:r>
    <w:t>Designed, constructed and assembled EMG data acquisition system for arm activities recognition</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="859"/>
      <w:tab w:val="left" w:pos="860"/>
    </w:tabs>
    <w:spacing w:line="244" w:lineRule="exact"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Converted time-domain data </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">of 200 gestures </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">into frequency domain using </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">fast </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>fourier</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> transform</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>to denoise signal</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="859"/>
      <w:tab w:val="left" w:pos="860"/>
    </w:tabs>
    <w:spacing w:line="244" w:lineRule="exact"/>
  </w:pPr>
  <w:r>
    <w:t>Classified different hand movements using support vector machines (</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>SVMs</w:t>
  </w:r>
  <w:r>
    <w:t>) with 82% accuracy</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="859"/>
      <w:tab w:val="left" w:pos="860"/>
    </w:tabs>
    <w:spacing w:line="244" w:lineRule="exact"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Improved accuracy by 3% training a </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:bCs/>
    </w:rPr>
    <w:t>neural network</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">providing insight for the </w:t>
  </w:r>
  <w:r>
    <w:t>medical rehabilitation system</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="right" w:pos="11080"/>
    </w:tabs>
    <w:ind w:left="144"/>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:color w:val="4472C4" w:themeColor="accent1"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="right" w:pos="11080"/>
    </w:tabs>
    <w:ind w:left="144"/>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:color w:val="4472C4" w:themeColor="accent1"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$startRng = $d.Range(0, 0)
$startRng.InsertXML($newContentXml)

# InsertXML silently drops a run-level <w:rStyle> (unlike a paragraph mark's
# rPr, which keeps it), so re-apply the "Strong" character style to the
# "Senior Researcher..." run explicitly.
$titleText = "Senior Researcher, Paul C. Lauterbur Lab at SIAT"
$titleRng = $d.Range(0, 0)
$titleRng.End = $titleRng.Start + $titleText.Length
$titleRng.Style = "Strong"

# ---------------------------------------------------------------------------
# Change 2: the Chinese "测试..." sentence was split across a page break into
# two runs; merge it back into a single run/sentence (removing the
# intervening lastRenderedPageBreak) by replacing the rendered text.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "测试，是如何测试的最后结果？测试的结果怎么样？为什么你做的东西有价值？这个环节能量化就量化。能量化的结果最有说服力。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "测试，是如何测试的最后结果？测试的结果怎么样？为什么你做的东西有价值？这个环节能量化就量化。能量化的结果最有说服力。",
    2)
